$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row at row 2; this pushes every existing data row
# down by one (old row N becomes row N+1), matching the target diff where
# the whole table shifts and a new record appears at the top.
$ws.Rows.Item(2).Insert()

# The freshly inserted row inherits formatting from the header row above it
# (bold font, borders, centered alignment). Reset it back to the plain,
# unstyled look used by every other data row.
$ws.Rows.Item(2).ClearFormats()

# Column D holds dates and needs the same custom date number format used by
# the rest of the table.
$ws.Cells.Item(2,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the new row's values.
$ws.Cells.Item(2,1).Value2  = 1
$ws.Cells.Item(2,2).Value2  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(2,3).Value2  = "Arica y Parinacota"
$ws.Cells.Item(2,4).Value2  = 44956
$ws.Cells.Item(2,5).Value2  = 15
$ws.Cells.Item(2,6).Value2  = 100112027
$ws.Cells.Item(2,7).Value2  = "Melón"
$ws.Cells.Item(2,8).Value2  = "Calameño"
$ws.Cells.Item(2,9).Value2  = "Super"
$ws.Cells.Item(2,10).Value2 = 100
$ws.Cells.Item(2,11).Value2 = 11000
$ws.Cells.Item(2,12).Value2 = 12000
$ws.Cells.Item(2,13).Value2 = 11500
$ws.Cells.Item(2,14).Value2 = "`$/caja 12 unidades"
$ws.Cells.Item(2,15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(2,16).Value2 = 958
$ws.Cells.Item(2,17).Value2 = 12
$ws.Cells.Item(2,18).Value2 = "Hortaliza"
